$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Index" header/table column to "i"
$ws.Range("A1").Value = "i"

# Re-base the index column from 1-based (1..502) to 0-based (0..501)
$n = 502
$arr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i,0] = $i
}
$ws.Range("A2:A503").Value = $arr

# Column A now fits a narrower header ("i" instead of "Index")
$ws.Range("A:A").ColumnWidth = 3.14
